# "fixed how datas are stored in stats"
#
# The "Log Stats" sheet stores, per row, a log-level label (column A) and a
# count (column B). This change:
#   - corrects the date shown in B1 ("2019-06-1" -> "2019-06-15")
#   - replaces the numeric placeholder labels in A4:A8 (0,1,2,3,4) with the
#     real log-level names (ERROR, INFO, DEBUG, FATAL, WARN)
#   - updates the counts in B2:B8 to the corrected figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 holds a date-shaped string ("2019-06-1" -> "2019-06-15"). Force the
# cell to text first so Excel stores the literal string instead of silently
# reinterpreting it as a date serial number, then drop the formatting again
# so the cell is left exactly as it started (General, no explicit style).
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2019-06-15"
$ws.Range("B1").ClearFormats()

# RandomLogs / main totals
$ws.Range("B2").Value = 25.0
$ws.Range("B3").Value = 25.0

# Log level rows: placeholder numeric codes -> real level names, plus the
# corrected counts for each level.
$ws.Range("A4").Value = "ERROR"
$ws.Range("B4").Value = 3.0

$ws.Range("A5").Value = "INFO"
$ws.Range("B5").Value = 3.0

$ws.Range("A6").Value = "DEBUG"
$ws.Range("B6").Value = 2.0

$ws.Range("A7").Value = "FATAL"
$ws.Range("B7").Value = 7.0

$ws.Range("A8").Value = "WARN"
$ws.Range("B8").Value = 10.0
